# Update cryptos list cell values (Price and Volume(1h) columns) to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.613.20"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.116.93"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").Value = "  +0.85%  "
$ws.Range("D5").Value = "'337.01"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("D7").Value = "'0.5252"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  +2.99%  "
$ws.Range("D9").Value = "'54.50"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").Value = "'0.09118"
$ws.Range("E10").Value = "  +2.22%  "
$ws.Range("D12").Value = "'24.45"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "2.118.98"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "'6.849"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("E15").Value = "  +5.77%  "
$ws.Range("D16").Value = "'0.00001175"
$ws.Range("E16").Value = "  +4.84%  "
$ws.Range("D17").Value = "'97.13"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "'1.012"
$ws.Range("D19").Value = "'0.06674"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "'19.40"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").Value = "'6.300"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").Value = "30.677.08"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("E24").Value = "  +5.05%  "
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("D26").Value = "2.357.34"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").Value = "'22.35"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("D28").Value = "'164.70"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").Value = "'2.552"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").Value = "'134.59"
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").Value = "'1.210"
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("D32").Value = "'0.1072"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "'1.646"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "'6.355"
$ws.Range("E34").Value = "  +3.46%  "
$ws.Range("D35").Value = "'3.945"
$ws.Range("E36").Value = "  +5.43%  "
$ws.Range("D37").Value = "'5.887"
$ws.Range("E37").Value = "  +7.78%  "
$ws.Range("D38").Value = "'0.02627"
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("D39").Value = "'0.06854"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").Value = "'0.2323"
$ws.Range("E40").Value = "  +3.18%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "'0.6891"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").Value = "'1.256"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "'14.80"
$ws.Range("E44").Value = "  +6.58%  "
$ws.Range("D45").Value = "'0.6487"
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("D46").Value = "'2.312"
$ws.Range("E46").Value = "  +5.47%  "
$ws.Range("D47").Value = "'0.00000000367"
$ws.Range("E47").Value = "  +21.43%  "
$ws.Range("D48").Value = "'3.692"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").Value = "'1.255"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("D50").Value = "'83.27"
$ws.Range("E50").Value = "  +2.13%  "
$ws.Range("D51").Value = "'1.192"
$ws.Range("E51").Value = "  -3.73%  "
